# Actualizacion de bitacora de respaldos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: backup to the repository was not actually done yet -> clear the date ---
$ws.Range("E16:F16").Value2 = $null

# --- Row 18: fix the "fecha en Repositorio" typo (was 42547, should mirror B18 = 42181) ---
$ws.Range("E18:F18").Value2 = 42181

# --- Row 19: new weekly backup log entry, formatted like the rows above it ---
$dateFmt = $ws.Range("B18").NumberFormat

$b19 = $ws.Range("B19")
$b19.Value2 = 42188
$b19.NumberFormat = $dateFmt
$b19.Borders.Item(7).LineStyle = 1
$b19.Borders.Item(8).LineStyle = 1
$b19.Borders.Item(9).LineStyle = 1
$b19.Borders.Item(10).LineStyle = 1

$cd19 = $ws.Range("C19:D19")
$cd19.Merge()
$cd19.Value2 = "Jovanny Zepeda"
$cd19.HorizontalAlignment = -4108
$cd19.Borders.Item(7).LineStyle = 1
$cd19.Borders.Item(8).LineStyle = 1
$cd19.Borders.Item(9).LineStyle = 1
$cd19.Borders.Item(10).LineStyle = 1

$ef19 = $ws.Range("E19:F19")
$ef19.Merge()
$ef19.Value2 = 42188
$ef19.NumberFormat = $dateFmt
$ef19.HorizontalAlignment = -4108
$ef19.Borders.Item(7).LineStyle = 1
$ef19.Borders.Item(8).LineStyle = 1
$ef19.Borders.Item(9).LineStyle = 1
$ef19.Borders.Item(10).LineStyle = 1

$g19 = $ws.Range("G19")
$g19.Value2 = "Realizada"
$g19.Borders.Item(7).LineStyle = 1
$g19.Borders.Item(8).LineStyle = 1
$g19.Borders.Item(9).LineStyle = 1
$g19.Borders.Item(10).LineStyle = 1

# --- Leave the selection where the user last edited ---
$ws.Range("E16:F16").Select() | Out-Null
